$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.555.28'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '1.640.28'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''212.50'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("E6").Value = '  +4.61%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''22.88'
$ws.Range("E8").Value = '  -4.52%  '
$ws.Range("D9").Value = '''0.256'
$ws.Range("E9").Value = '  -1.70%  '
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("D12").Value = '1.873.77'
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("D13").Value = '1.652.15'
$ws.Range("E13").Value = '  -0.09%  '
$ws.Range("E14").Value = '  -0.99%  '
$ws.Range("E15").Value = '  -1.24%  '
$ws.Range("D16").Value = '''63.94'
$ws.Range("E16").Value = '  -2.61%  '
$ws.Range("D17").Value = '27.528.32'
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").Value = '''227.50'
$ws.Range("E18").Value = '  -1.73%  '
$ws.Range("E19").Value = '  -0.30%  '
$ws.Range("D20").Value = '''7.60'
$ws.Range("E20").Value = '  +1.57%  '
$ws.Range("D22").Value = '''4.29'
$ws.Range("E22").Value = '  -1.81%  '
$ws.Range("D23").Value = '''10.00'
$ws.Range("E23").Value = '  +7.59%  '
$ws.Range("D24").Value = '''1.95'
$ws.Range("E24").Value = '  -3.74%  '
$ws.Range("D25").Value = '''149.21'
$ws.Range("E25").Value = '  +1.61%  '
$ws.Range("D26").Value = '''6.95'
$ws.Range("E26").Value = '  -3.13%  '
$ws.Range("E27").Value = '  +1.32%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").Value = '''15.58'
$ws.Range("E29").Value = '  -1.84%  '
$ws.Range("E31").Value = '  -2.40%  '
$ws.Range("D32").Value = '''3.29'
$ws.Range("E32").Value = '  -0.33%  '
$ws.Range("E33").Value = '  +1.71%  '
$ws.Range("D34").Value = '1.426.31'
$ws.Range("E34").Value = '  -2.26%  '
$ws.Range("E35").Value = '  +2.09%  '
$ws.Range("E36").Value = '  -1.96%  '
$ws.Range("D37").Value = '''0.572'
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("D38").Value = '''0.874'
$ws.Range("E38").Value = '  -3.76%  '
$ws.Range("E39").Value = '  -1.26%  '
$ws.Range("D40").Value = '''0.900'
$ws.Range("E40").Value = '  +14.96%  '
$ws.Range("E41").Value = '  -2.17%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("B43").Value = 'mCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D43").Value = '''2.46'
$ws.Range("E43").Value = '  -0.97%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '''5.50'
$ws.Range("E44").Value = '  +0.91%  '
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").Value = '''2.25'
$ws.Range("E45").Value = '  +1.39%  '
$ws.Range("D46").Value = '''65.08'
$ws.Range("E46").Value = '  -0.35%  '
$ws.Range("D47").Value = '1.783.01'
$ws.Range("E47").Value = '  -0.68%  '
$ws.Range("E48").Value = '  -2.36%  '
$ws.Range("D49").Value = '''86.25'
$ws.Range("E49").Value = '  -2.26%  '
$ws.Range("E50").Value = '  +1.25%  '
$ws.Range("D51").Value = '''0.0983'
$ws.Range("E51").Value = '  -2.64%  '
